# Auto-generated edit script applying scheduled-runner price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1552.3334
$ws.Range("I8").Value = 1298.6364
$ws.Range("J8").Value = 2250
$ws.Range("K8").Value = 3895.9092
$ws.Range("L8").Value = 6750
$ws.Range("M8").Value = -3756.9092
$ws.Range("N8").Value = -7028
$ws.Range("H18").Value = 1033.3334
$ws.Range("I18").Value = 800.25
$ws.Range("K18").Value = 800.25
$ws.Range("M18").Value = -516.25
$ws.Range("H32").Value = 6098.75
$ws.Range("J32").Value = 5134.778
$ws.Range("L32").Value = 5134.778
$ws.Range("N32").Value = -5786.778
$ws.Range("H76").Value = 9099.799999999999
$ws.Range("J76").Value = 9222
$ws.Range("L76").Value = 9222
$ws.Range("N76").Value = -9852
$ws.Range("H79").Value = 9099.799999999999
$ws.Range("J79").Value = 9222
$ws.Range("L79").Value = 9222
$ws.Range("N79").Value = -11406
$ws.Range("H134").Value = 66081.5
$ws.Range("J134").Value = 66081.5
$ws.Range("L134").Value = 66081.5
$ws.Range("N134").Value = -76221.5
$ws.Range("H141").Value = 3282.4
$ws.Range("I141").Value = 2728
$ws.Range("K141").Value = 8184
$ws.Range("M141").Value = -3004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 13342.429
$ws.Range("I3").Value = 1975.5
$ws.Range("J3").Value = 28498.334
$ws.Range("K3").Value = 1975.5
$ws.Range("L3").Value = 28498.334
$ws.Range("M3").Value = -1860.5
$ws.Range("N3").Value = -28728.334
$ws.Range("H12").Value = 1999.3334
$ws.Range("J12").Value = 998
$ws.Range("L12").Value = 998
$ws.Range("N12").Value = -1344
$ws.Range("H13").Value = 3339168
$ws.Range("J13").Value = 8752
$ws.Range("L13").Value = 8752
$ws.Range("N13").Value = -9040
$ws.Range("H17").Value = 3800
$ws.Range("I17").Value = 3800
$ws.Range("J17").Value = 3800
$ws.Range("K17").Value = 3800
$ws.Range("L17").Value = 3800
$ws.Range("M17").Value = -3627
$ws.Range("N17").Value = -4146
$ws.Range("H18").Value = 999
$ws.Range("J18").Value = 999
$ws.Range("L18").Value = 999
$ws.Range("H19").Value = 5027.75
$ws.Range("I19").Value = 5050.5
$ws.Range("K19").Value = 5050.5
$ws.Range("M19").Value = -4821.5
$ws.Range("H74").Value = 13337421
$ws.Range("I74").Value = 22224564
$ws.Range("J74").Value = 6706.7
$ws.Range("K74").Value = 22224564
$ws.Range("L74").Value = 6706.7
$ws.Range("M74").Value = -22223690
$ws.Range("N74").Value = -8454.700000000001
$ws.Range("H77").Value = 13337421
$ws.Range("I77").Value = 22224564
$ws.Range("J77").Value = 6706.7
$ws.Range("K77").Value = 111122820
$ws.Range("L77").Value = 33533.5
$ws.Range("M77").Value = -111118452
$ws.Range("N77").Value = -42269.5
$ws.Range("N18").Value = -1643

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 42925
$ws.Range("J60").Value = 42925
$ws.Range("L60").Value = 42925
$ws.Range("N60").Value = -44123
$ws.Range("H86").Value = 4732.643
$ws.Range("J86").Value = 10101.75
$ws.Range("L86").Value = 10101.75
$ws.Range("N86").Value = -12347.75
$ws.Range("H89").Value = 4732.643
$ws.Range("J89").Value = 10101.75
$ws.Range("L89").Value = 50508.75
$ws.Range("N89").Value = -61740.75
$ws.Range("H94").Value = 1390.2273
$ws.Range("I94").Value = 1446.1052
$ws.Range("K94").Value = 1446.1052
$ws.Range("M94").Value = -995.1052
$ws.Range("H99").Value = 1566.8667
$ws.Range("I99").Value = 1173.4546
$ws.Range("K99").Value = 1173.4546
$ws.Range("M99").Value = 324.5454
$ws.Range("H109").Value = 72870.5
$ws.Range("J109").Value = 72870.5
$ws.Range("L109").Value = 72870.5
$ws.Range("N109").Value = -75644.5
$ws.Range("H134").Value = 2302.4666
$ws.Range("I134").Value = 2302.4666
$ws.Range("K134").Value = 6907.399800000001
$ws.Range("M134").Value = -4372.399800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31360.82
$ws.Range("I31").Value = 3838.1333
$ws.Range("K31").Value = 3838.1333
$ws.Range("M31").Value = -3543.1333
$ws.Range("H34").Value = 31360.82
$ws.Range("I34").Value = 3838.1333
$ws.Range("K34").Value = 3838.1333
$ws.Range("M34").Value = -3636.1333
$ws.Range("H132").Value = 3346.7144
$ws.Range("I132").Value = 2589.3044
$ws.Range("K132").Value = 7767.9132
$ws.Range("M132").Value = -5237.9132
$ws.Range("H134").Value = 3818.077
$ws.Range("I134").Value = 2847
$ws.Range("K134").Value = 8541
$ws.Range("M134").Value = -6006

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 7004.5
$ws.Range("I119").Value = 1782.2222
$ws.Range("K119").Value = 5346.6666
$ws.Range("M119").Value = -508.6665999999996
$ws.Range("H131").Value = 7793109
$ws.Range("I131").Value = 15626154
$ws.Range("J131").Value = 5894188.5
$ws.Range("K131").Value = 46878462
$ws.Range("L131").Value = 17682565.5
$ws.Range("M131").Value = -46873422
$ws.Range("N131").Value = -17692645.5
$ws.Range("H132").Value = 4667.4287
$ws.Range("I132").Value = 3791.8572
$ws.Range("J132").Value = 5543
$ws.Range("K132").Value = 34126.7148
$ws.Range("L132").Value = 49887
$ws.Range("M132").Value = -31596.7148
$ws.Range("N132").Value = -54947

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12841.866
$ws.Range("I70").Value = 9848.556
$ws.Range("K70").Value = 9848.556
$ws.Range("M70").Value = -9578.556
$ws.Range("H73").Value = 12841.866
$ws.Range("I73").Value = 9848.556
$ws.Range("K73").Value = 9848.556
$ws.Range("M73").Value = -8912.556
$ws.Range("H80").Value = 337124.06
$ws.Range("I80").Value = 716651.4399999999
$ws.Range("J80").Value = 5037.625
$ws.Range("K80").Value = 716651.4399999999
$ws.Range("L80").Value = 5037.625
$ws.Range("M80").Value = -715653.4399999999
$ws.Range("N80").Value = -7033.625
$ws.Range("H83").Value = 337124.06
$ws.Range("I83").Value = 716651.4399999999
$ws.Range("J83").Value = 5037.625
$ws.Range("K83").Value = 3583257.2
$ws.Range("L83").Value = 25188.125
$ws.Range("M83").Value = -3578265.2
$ws.Range("N83").Value = -35172.125
$ws.Range("H132").Value = 3970.923
$ws.Range("I132").Value = 2987
$ws.Range("K132").Value = 8961
$ws.Range("M132").Value = -6431

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2106.5833
$ws.Range("I93").Value = 2106.5833
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2106.5833
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -858.5832999999998
$ws.Range("H136").Value = 4586.4346
$ws.Range("I136").Value = 3568.2903
$ws.Range("J136").Value = 6690.6
$ws.Range("K136").Value = 10704.8709
$ws.Range("L136").Value = 20071.8
$ws.Range("M136").Value = -8154.8709
$ws.Range("N136").Value = -25171.8
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9930.5
$ws.Range("I132").Value = 4608.5
$ws.Range("J132").Value = 15252.5
$ws.Range("K132").Value = 13825.5
$ws.Range("L132").Value = 45757.5
$ws.Range("M132").Value = -11295.5
$ws.Range("N132").Value = -50817.5
$ws.Range("H136").Value = 4650.6665
$ws.Range("I136").Value = 2705.0908
$ws.Range("J136").Value = 10001
$ws.Range("K136").Value = 8115.2724
$ws.Range("L136").Value = 30003
$ws.Range("M136").Value = -5565.2724
$ws.Range("N136").Value = -35103
